$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values (row 3 = Luty, row 9 = Sierpien, row 10 = Wrzesien, row 12 = Listopad)
$ws.Range("B3").Value = 168
$ws.Range("C3").Value = 21

$ws.Range("B9").Value = 168
$ws.Range("C9").Value = 21
$ws.Range("D9").Value = 10

$ws.Range("B10").Value = 168
$ws.Range("C10").Value = 21
$ws.Range("D10").Value = 9

$ws.Range("B12").Value = 152
$ws.Range("C12").Value = 19
$ws.Range("D12").Value = 11

# Change selection
$ws.Range("H7").Select()
